# Applies the "Global model tune after test loop fix" update:
# - adds a second (F:H) MASE stats block next to the existing (B:C) block
#   on rows 24-43 of the active sheet ("Sheet1")
# - updates the sheet view selection / scroll position
# - updates the workbook view scroll position
# - sets the page setup (paper size / orientation) for the sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New header label (row 24, column F) ---
$ws.Range("F24").Value = "Global model tune after test loop fix"

# --- New values for the second MASE summary block ---
$ws.Range("F25").Value = 1.0960000000000001
$ws.Range("G25").Value = "Column1"
$ws.Range("H25").Value = ""

$ws.Range("F26").Value = 0.78256000000000003

$ws.Range("F27").Value = 1.2548600000000001
$ws.Range("G27").Value = "Mean"
$ws.Range("H27").Value = 1.440548947368421

$ws.Range("F28").Value = 0.71653999999999995
$ws.Range("G28").Value = "Standard Error"
$ws.Range("H28").Value = 0.34545852209356903

$ws.Range("F29").Value = 1.1507000000000001
$ws.Range("G29").Value = "Median"
$ws.Range("H29").Value = 1.07944

$ws.Range("F30").Value = 0.71194000000000002
$ws.Range("G30").Value = "Mode"
$ws.Range("H30").Value = "#N/A"

$ws.Range("F31").Value = 0.72307999999999995
$ws.Range("G31").Value = "Standard Deviation"
$ws.Range("H31").Value = 1.5058187869907806

$ws.Range("F32").Value = 1.4930000000000001
$ws.Range("G32").Value = "Sample Variance"
$ws.Range("H32").Value = 2.2674902192543862

$ws.Range("F33").Value = 1.0366
$ws.Range("G33").Value = "Kurtosis"
$ws.Range("H33").Value = 15.146573058098053

$ws.Range("F34").Value = 1.67167
$ws.Range("G34").Value = "Skewness"
$ws.Range("H34").Value = 3.7517685480969303

$ws.Range("F35").Value = 0.5736
$ws.Range("G35").Value = "Range"
$ws.Range("H35").Value = 6.7971900000000005

$ws.Range("F36").Value = 2.4283700000000001
$ws.Range("G36").Value = "Minimum"
$ws.Range("H36").Value = 0.5736

$ws.Range("F37").Value = 1.50444
$ws.Range("G37").Value = "Maximum"
$ws.Range("H37").Value = 7.3707900000000004

$ws.Range("F38").Value = 1.38266
$ws.Range("G38").Value = "Sum"
$ws.Range("H38").Value = 27.370429999999999

$ws.Range("F39").Value = 0.57755999999999996
$ws.Range("G39").Value = "Count"
$ws.Range("H39").Value = 19

$ws.Range("F40").Value = 0.83343
$ws.Range("H40").Value = 0

$ws.Range("F41").Value = 1.07944

$ws.Range("F42").Value = 7.3707900000000004

$ws.Range("F43").Value = 0.98319000000000001

# --- Copy the existing "Column1 / Mean / Standard Error / ..." stat block
#     (columns B:C, rows 25-39) formatting onto the new columns (G:H), now
#     that the values/text are already in place ---
$ws.Range("B25:C39").Copy() | Out-Null
$ws.Range("G25:G39").PasteSpecial(-4122) | Out-Null
$ws.Range("C25:C39").Copy() | Out-Null
$ws.Range("H25:H39").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- View / selection updates ---
$ws.Range("F25").Select() | Out-Null

$wb.Windows.Item(1).Top = 375

# --- Page setup ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
